$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 25 (Bíobío / Vega Monumental Concepción weekly log is
# prepended with a new daily observation; every existing row from 25..102
# shifts down to 26..103).
$ws.Rows.Item(25).Insert()

# Populate the newly inserted row 25 with the new observation.
$ws.Range("A25").Value = 11
$ws.Range("B25").Value = "Vega Monumental Concepción"
$ws.Range("C25").Value = "Bíobío"
$ws.Range("D25").Value = 44623
$ws.Range("E25").Value = 8
$ws.Range("F25").Value = 100112032
$ws.Range("G25").Value = "Zapallo italiano"
$ws.Range("H25").Value = "Sin especificar"
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 220
$ws.Range("K25").Value = 10000
$ws.Range("L25").Value = 11000
$ws.Range("M25").Value = 10455
$ws.Range("N25").Value = "$/caja 60 unidades"
$ws.Range("O25").Value = "Región de Arica y Parinacota"
$ws.Range("P25").Value = 174
$ws.Range("Q25").Value = 60
$ws.Range("R25").Value = "Hortaliza"
